$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.254.24'
$ws.Range("E2").Value = '  -0.69%  '
$ws.Range("D3").Value = '1.875.15'
$ws.Range("E3").Value = '  -1.87%  '
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '234.78'
$ws.Range("E5").Value = '  -1.92%  '
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").Value = '0.4699'
$ws.Range("E7").Value = '  -1.76%  '
$ws.Range("D8").Value = '0.2837'
$ws.Range("E8").Value = '  -0.20%  '
$ws.Range("D9").Value = '0.06605'
$ws.Range("E9").Value = '  -1.42%  '
$ws.Range("D10").Value = '20.31'
$ws.Range("E10").Value = '  +7.89%  '
$ws.Range("D11").Value = '0.07775'
$ws.Range("E11").Value = '  +0.91%  '
$ws.Range("D12").Value = '97.72'
$ws.Range("E12").Value = '  -4.20%  '
$ws.Range("D13").Value = '1.880.65'
$ws.Range("E13").Value = '  -2.02%  '
$ws.Range("D14").Value = '5.085'
$ws.Range("E14").Value = '  -2.37%  '
$ws.Range("D15").Value = '0.6735'
$ws.Range("E15").Value = '  +0.46%  '
$ws.Range("D16").Value = '287.52'
$ws.Range("E16").Value = '  +8.28%  '
$ws.Range("D17").Value = '30.284.42'
$ws.Range("E17").Value = '  -0.69%  '
$ws.Range("E18").Value = '  -0.02%  '
$ws.Range("E19").Value = '  -0.53%  '
$ws.Range("D20").Value = '2.129.29'
$ws.Range("E20").Value = '  -1.36%  '
$ws.Range("D21").Value = '5.384'
$ws.Range("E21").Value = '  -0.58%  '
$ws.Range("D22").Value = '0.000007279'
$ws.Range("E22").Value = '  -2.59%  '
$ws.Range("D23").Value = '0.9999'
$ws.Range("E23").Value = '  -0.16%  '
$ws.Range("D24").Value = '6.177'
$ws.Range("E24").Value = '  -1.89%  '
$ws.Range("D25").Value = '9.386'
$ws.Range("E25").Value = '  +0.13%  '
$ws.Range("D26").Value = '167.45'
$ws.Range("E26").Value = '  +0.36%  '
$ws.Range("E27").Value = '  -0.15%  '
$ws.Range("D28").Value = '1.975'
$ws.Range("E28").Value = '  -4.33%  '
$ws.Range("D29").Value = '1.370'
$ws.Range("E29").Value = '  -1.11%  '
$ws.Range("D30").Value = '0.09674'
$ws.Range("E30").Value = '  -3.38%  '
$ws.Range("D31").Value = '4.393'
$ws.Range("E31").Value = '  -5.32%  '
$ws.Range("D32").Value = '1.467'
$ws.Range("E32").Value = '  -3.14%  '
$ws.Range("D33").Value = '4.120'
$ws.Range("E33").Value = '  -2.44%  '
$ws.Range("D34").Value = '0.04690'
$ws.Range("E34").Value = '  -0.82%  '
$ws.Range("D36").Value = '1.093'
$ws.Range("E36").Value = '  -1.16%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '2.717'
$ws.Range("E37").Value = '  -0.14%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.01871'
$ws.Range("E38").Value = '  -2.18%  '
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = '6.429'
$ws.Range("E39").Value = '  +2.50%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '2.526'
$ws.Range("E40").Value = '  -3.29%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = '72.04'
$ws.Range("E41").Value = '  -4.00%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '1.957'
$ws.Range("E42").Value = '  -0.56%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '0.8613'
$ws.Range("E43").Value = '  +0.16%  '
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = '0.9998'
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '103.13'
$ws.Range("E45").Value = '  -1.77%  '
$ws.Range("B46").Value = 'TheSandbox'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D46").Value = '0.4196'
$ws.Range("E46").Value = '  -1.66%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '989.51'
$ws.Range("E47").Value = '  +7.23%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '7.219'
$ws.Range("E48").Value = '  -2.51%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '9.151'
$ws.Range("E49").Value = '  +3.89%  '
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '34.01'
$ws.Range("E50").Value = '  -2.26%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = '0.1150'
$ws.Range("E51").Value = '  -4.24%  '
